# Generate Report for Handback
# This mirrors a localization pipeline run that recorded a handback
# transform failure for the 6507b6a1-... file: the status flips from
# "Ready for handoff" to "Handback transform failed" on the Overview
# sheet as well as on the per-locale (zh-cn / de-de) detail sheets, and
# the corresponding "Error Detail" cell on each locale sheet is filled
# in with the mismatch explanation. The Error Detail column is also
# widened so the message is readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$zhError = "Handback file name: 4bdjkc2p.xxb is different with handoff file name: 6507b6a1-e03d-4865-a9b5-f79a4e1dd8c3.7cbcefa6815ce2c7415173d4497dd2a41ca1d5ff.zh-cn."
$deError  = "Handback file name: 4bdjkc2p.xxb is different with handoff file name: 6507b6a1-e03d-4865-a9b5-f79a4e1dd8c3.7cbcefa6815ce2c7415173d4497dd2a41ca1d5ff.de-de."

# --- Overview sheet: row 3 is the 6507b6a1 file; E3/F3 hold the status
#     shown for zh-cn / de-de respectively.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn detail sheet: row 3 is the 6507b6a1 file.
# (ColumnWidth is specified in characters; the stored OOXML column
#  width is in pixel-derived units, so 39.1666... characters is what
#  round-trips to a stored width of exactly 40.)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de detail sheet: row 3 is the 6507b6a1 file.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
